# "add VS3 model 1"
#
# Inserts a new worksheet "act_model_1_VS_3" right after "act model 1"
# (so it becomes the 3rd tab). It holds the VS_B / VS_L / VS_R
# "act_model_1" contrast columns (I, P, W on the main dashboard sheet,
# rows 2:45) as plain literal values - mirroring the existing
# "act_model_3_VS_3" sheet that already sits next to it.
#
# The new sheet becomes the active tab, and the selection on the main
# dashboard sheet moves from the old "model 3" VS_R column (Y) to the
# new "model 1" VS_R column (W) to reflect the newly-added model.

$wb = $excel.ActiveWorkbook

# --- locate the anchor sheet ("act model 1") and the main dashboard sheet
$mainSheet = $wb.Worksheets.Item(1)
$actModel1 = $wb.Worksheets.Item("act model 1")

# --- insert the new sheet right after "act model 1" and rename it
$newSheet = $wb.Worksheets.Add($null, $actModel1)
$newSheet.Name = "act_model_1_VS_3"

# --- copy VS_B (col I=9), VS_L (col P=16), VS_R (col W=23) act_model_1
#     values from the dashboard sheet (rows 2:45) into A:C of the new sheet
for ($r = 2; $r -le 45; $r++) {
    $destRow = $r - 1
    $newSheet.Cells.Item($destRow, 1).Value = $mainSheet.Cells.Item($r, 9).Value2
    $newSheet.Cells.Item($destRow, 2).Value = $mainSheet.Cells.Item($r, 16).Value2
    $newSheet.Cells.Item($destRow, 3).Value = $mainSheet.Cells.Item($r, 23).Value2
}

# --- move the selection on the main dashboard sheet onto the new
#     act_model_1 VS_R column (was VS_R act_model_3 / column Y)
$mainSheet.Activate() | Out-Null
$mainSheet.Range("W2:W45").Select() | Out-Null

# --- select the whole data block on the new sheet and make it the active tab
$newSheet.Activate() | Out-Null
$newSheet.Range("A1:C44").Select() | Out-Null
